# PlayerPerformance_4099.xlsx update:
#  1. Insert a new "Player Info" worksheet as the first sheet, with player
#     metadata (ID / NAME / BATTING_HAND / BOWL_STYLE).
#  2. In the existing "ODI Batting" sheet, rename the MATCH_CARD_LINK column
#     to MATCH_CODE and replace each full scorecard URL with the bare
#     numeric match code extracted from it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "Player Info" sheet (becomes sheet index 1 / first tab).
# ---------------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Match the bold/bordered/centered-top header style already used elsewhere
# in this workbook.
$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Data row. ID is kept as text (not a number) to match the source data feed.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "4099"
$playerInfo.Range("A2").ClearFormats()

$playerInfo.Range("B2").Value = "Sanju Viswanath Samson"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

# ---------------------------------------------------------------------------
# 2. "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE.
# ---------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"

$matchCodes = @("4485", "4621", "4623", "4624", "4637", "4640", "4643", "4656", "4657", "4658", "4669")

for ($i = 0; $i -lt $matchCodes.Length; $i++) {
    $row = $i + 2
    $cell = $odiBatting.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $matchCodes[$i]
    $cell.ClearFormats()
}
